# Insert a new record (row 446) into the Coliflor / Macroferia Regional de Talca
# daily-price log. This pushes every existing record from the old row 446
# onward down by one row (old row 445 stays put; old row 545 becomes the new
# row 546), and the freshly inserted row 446 is populated with a new
# observation (03/12/2023 -> Excel serial 45275) while the constant columns
# (market, region, product, unit, origin, classification, etc.) are copied
# from the surrounding rows since every record in this sheet shares them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 446:545 down to 447:546, opening up a blank row 446.
$ws.Rows.Item(446).Insert()

# Populate the newly inserted row 446 with the new observation.
$ws.Range("A446").Value = 5
$ws.Range("B446").Value = "Macroferia Regional de Talca"
$ws.Range("C446").Value = "Maule"
$ws.Range("D446").Value = 45275
$ws.Range("E446").Value = 7
$ws.Range("F446").Value = 100112008
$ws.Range("G446").Value = "Coliflor"
$ws.Range("H446").Value = "Sin especificar"
$ws.Range("I446").Value = "Primera"
$ws.Range("J446").Value = 3000
$ws.Range("K446").Value = 800
$ws.Range("L446").Value = 800
$ws.Range("M446").Value = 800
$ws.Range("N446").Value = '$/unidad'
$ws.Range("O446").Value = "Región del Maule"
$ws.Range("P446").Value = 800
$ws.Range("Q446").Value = 1
$ws.Range("R446").Value = "Hortaliza"
